$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C55").Value = "speedlimit"
$ws.Range("C80").Value = "speedlimit"
$ws.Range("C81").Value = "speedlimit"
$ws.Range("C83").Value = "speedlimit"
$ws.Range("C84").Value = "other"
$ws.Range("C85").Value = "speedlimit"
$ws.Range("C91").Value = "other"
$ws.Range("C94").Value = "speedlimit"
$ws.Range("C102").Value = "other"
$ws.Range("C104").Value = "speedlimit"
$ws.Range("C105").Value = "speedlimit"
$ws.Range("C135").Value = "speedlimit"
$ws.Range("C145").Value = "speedlimit"
$ws.Range("C172").Value = "other"
$ws.Range("C190").Value = "speedlimit"
$ws.Range("C197").Value = "other"
$ws.Range("C198").Value = "other"
$ws.Range("C210").Value = "other"
$ws.Range("C222").Value = "speedlimit"
$ws.Range("C233").Value = "other"
$ws.Range("C248").Value = "speedlimit"
$ws.Range("C252").Value = "other"
$ws.Range("C261").Value = "speedlimit"
$ws.Range("C268").Value = "speedlimit"
$ws.Range("C273").Value = "other"
$ws.Range("C274").Value = "other"
$ws.Range("C277").Value = "speedlimit"
$ws.Range("C312").Value = "speedlimit"
$ws.Range("C313").Value = "speedlimit"
$ws.Range("C327").Value = "speedlimit"
$ws.Range("C356").Value = "other"
$ws.Range("C363").Value = "speedlimit"
$ws.Range("C411").Value = "other"
$ws.Range("C419").Value = "speedlimit"
$ws.Range("C420").Value = "speedlimit"
$ws.Range("C430").Value = "other"
$ws.Range("C442").Value = "speedlimit"
$ws.Range("C447").Value = "other"
$ws.Range("C448").Value = "other"
$ws.Range("C462").Value = "speedlimit"
$ws.Range("C467").Value = "speedlimit"
$ws.Range("C483").Value = "other"
$ws.Range("C490").Value = "other"
$ws.Range("C492").Value = "speedlimit"
$ws.Range("C504").Value = "other"
$ws.Range("C505").Value = "other"
$ws.Range("C512").Value = "speedlimit"
$ws.Range("C516").Value = "speedlimit"
$ws.Range("C524").Value = "speedlimit"
